$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Range("D2")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "42.003.77"
$cD.Style = $styleD
$ws.Range("E2").Value = "  -1.68%  "

$cD = $ws.Range("D3")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.290.62"
$cD.Style = $styleD
$ws.Range("E3").Value = "  -2.52%  "

$ws.Range("E4").Value = "  -0.01%  "

$cD = $ws.Range("D5")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "313.81"
$cD.Style = $styleD
$ws.Range("E5").Value = "  -2.05%  "

$cD = $ws.Range("D6")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "103.76"
$cD.Style = $styleD
$ws.Range("E6").Value = "  -1.96%  "

$ws.Range("E7").Value = "  -2.14%  "

$ws.Range("E8").Value = "  -0.06%  "

$cD = $ws.Range("D9")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.607"
$cD.Style = $styleD
$ws.Range("E9").Value = "  -1.55%  "

$cD = $ws.Range("D10")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "39.78"
$cD.Style = $styleD
$ws.Range("E10").Value = "  -2.84%  "

$cD = $ws.Range("D11")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.0909"
$cD.Style = $styleD
$ws.Range("E11").Value = "  -1.61%  "

$cD = $ws.Range("D12")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "8.40"
$cD.Style = $styleD
$ws.Range("E12").Value = "  -0.01%  "

$cD = $ws.Range("D13")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.106"
$cD.Style = $styleD
$ws.Range("E13").Value = "  +0.57%  "

$cD = $ws.Range("D14")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.969"
$cD.Style = $styleD
$ws.Range("E14").Value = "  -2.41%  "

$cD = $ws.Range("D15")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "15.27"
$cD.Style = $styleD
$ws.Range("E15").Value = "  -4.59%  "

$cD = $ws.Range("D16")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.637.16"
$cD.Style = $styleD
$ws.Range("E16").Value = "  -2.47%  "

$cD = $ws.Range("D17")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.307.68"
$cD.Style = $styleD
$ws.Range("E17").Value = "  -1.68%  "

$cD = $ws.Range("D18")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "41.932.37"
$cD.Style = $styleD
$ws.Range("E18").Value = "  -1.71%  "

$cD = $ws.Range("D19")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "7.55"
$cD.Style = $styleD
$ws.Range("E19").Value = "  -1.96%  "

$ws.Range("E20").Value = "  -1.12%  "

$cD = $ws.Range("D21")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "72.60"
$cD.Style = $styleD
$ws.Range("E21").Value = "  -5.92%  "

$cD = $ws.Range("D22")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "3.52"
$cD.Style = $styleD
$ws.Range("E22").Value = "  -1.27%  "

$cD = $ws.Range("D23")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "257.03"
$cD.Style = $styleD
$ws.Range("E23").Value = "  -1.17%  "

$cD = $ws.Range("D24")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.30"
$cD.Style = $styleD
$ws.Range("E24").Value = "  -0.91%  "

$cD = $ws.Range("D25")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "9.78"
$cD.Style = $styleD
$ws.Range("E25").Value = "  +2.33%  "

$cD = $ws.Range("D26")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.00"
$cD.Style = $styleD
$ws.Range("E26").Value = "  +0.40%  "

$cD = $ws.Range("D27")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "10.93"
$cD.Style = $styleD
$ws.Range("E27").Value = "  -4.40%  "

$ws.Range("E28").Value = "  +2.49%  "

$cD = $ws.Range("D29")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "22.68"
$cD.Style = $styleD
$ws.Range("E29").Value = "  -2.66%  "

$cD = $ws.Range("D30")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "35.51"
$cD.Style = $styleD
$ws.Range("E30").Value = "  -1.79%  "

$cD = $ws.Range("D31")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "163.33"
$cD.Style = $styleD
$ws.Range("E31").Value = "  -6.39%  "

$cD = $ws.Range("D32")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.0883"
$cD.Style = $styleD
$ws.Range("E32").Value = "  -0.87%  "

$ws.Range("E33").Value = "  -2.68%  "

$cD = $ws.Range("D34")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "5.84"
$cD.Style = $styleD
$ws.Range("E34").Value = "  -3.70%  "

$ws.Range("E35").Value = "  -0.81%  "

$ws.Range("E36").Value = "  +4.86%  "

$cD = $ws.Range("D37")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "4.60"
$cD.Style = $styleD
$ws.Range("E37").Value = "  -0.56%  "

$ws.Range("E38").Value = "  +8.53%  "

$ws.Range("E39").Value = "  -2.19%  "

$cD = $ws.Range("D40")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "3.60"
$cD.Style = $styleD
$ws.Range("E40").Value = "  -5.22%  "

$cD = $ws.Range("D41")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "99.06"
$cD.Style = $styleD
$ws.Range("E41").Value = "  +18.35%  "

$ws.Range("E42").Value = "  +0.90%  "

$cD = $ws.Range("D43")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "70.31"
$cD.Style = $styleD
$ws.Range("E43").Value = "  -0.37%  "

$cD = $ws.Range("D44")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.226"
$cD.Style = $styleD
$ws.Range("E44").Value = "  -2.53%  "

$cD = $ws.Range("D46")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "12.04"
$cD.Style = $styleD
$ws.Range("E46").Value = "  +1.20%  "

$cD = $ws.Range("D47")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "112.94"
$cD.Style = $styleD
$ws.Range("E47").Value = "  -1.20%  "

$cD = $ws.Range("D48")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "77.82"
$cD.Style = $styleD
$ws.Range("E48").Value = "  +6.67%  "

$cD = $ws.Range("D49")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "9.11"
$cD.Style = $styleD
$ws.Range("E49").Value = "  -0.96%  "

$cD = $ws.Range("D50")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "5.29"
$cD.Style = $styleD
$ws.Range("E50").Value = "  -3.97%  "

$cD = $ws.Range("D51")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.29"
$cD.Style = $styleD
$ws.Range("E51").Value = "  +2.59%  "
